$d = $word.ActiveDocument

# --- Introduction body paragraph (SVM -> Newton's laws) ---
$d.Paragraphs(2).Range.Text = @"
Newton’s laws of motion, formulated by Sir Isaac Newton in 1687, are foundational principles in physics that describe the relationship between a body and the forces acting upon it. These three laws—often referred to as the law of inertia, the law of acceleration, and the law of action and reaction—form the cornerstone of classical mechanics. They explain how objects move, respond to forces, and interact with one another. Despite being over three centuries old, Newton’s laws remain essential in understanding and predicting the behavior of physical systems, from everyday phenomena to advanced technological applications.
"@

# --- Heading2: "Real-World Applications of SVM" -> "Real World Applications of the Topics" ---
# Rename bookmark id=21: real-world-applications-of-svm -> real-world-applications-of-the-topics
$bm = $d.Bookmarks("real-world-applications-of-svm")
$rng = $bm.Range
$bm.Delete()
$d.Bookmarks.Add("real-world-applications-of-the-topics", $rng)
$d.Paragraphs(3).Range.Text = "Real World Applications of the Topics"

# --- Heading3: "Healthcare" -> "Aerospace Engineering" (bookmark id=22) ---
$bm = $d.Bookmarks("healthcare")
$rng = $bm.Range
$bm.Delete()
$d.Bookmarks.Add("aerospace-engineering", $rng)
$d.Paragraphs(4).Range.Text = "Aerospace Engineering"

# --- Body paragraph under "Aerospace Engineering" ---
$d.Paragraphs(5).Range.Text = @"
Newton’s laws are fundamental to aerospace engineering, where they are used to design and control aircraft and spacecraft. The first law of motion (inertia) is crucial for understanding how vehicles maintain their trajectory in space. The second law (force and acceleration) helps calculate the thrust and fuel requirements for rockets, while the third law (action and reaction) explains how rockets propel themselves forward by expelling gases. For instance, SpaceX relies on Newton’s laws to optimize rocket launches and ensure precise orbital insertions.
"@

# --- Heading3: "Stock Market Forecasting" -> "Robotics" (bookmark id=23) ---
$bm = $d.Bookmarks("stock-market-forecasting")
$rng = $bm.Range
$bm.Delete()
$d.Bookmarks.Add("robotics", $rng)
$d.Paragraphs(6).Range.Text = "Robotics"

# --- Body paragraph under "Robotics" ---
$d.Paragraphs(7).Range.Text = @"
Robots operate based on Newton’s laws of motion, which guide their movement and interaction with their environment. The first law ensures robots maintain their motion until acted upon by external forces, while the second law helps calculate the torque and force required for precise mechanical movements. The third law is essential for understanding the forces exchanged between a robot and its surroundings. Modern robots, such as those in manufacturing or delivery systems, rely heavily on these principles to perform tasks efficiently and safely.
"@

# --- Heading3: "Predictive Maintenance" -> "Automotive Safety" (bookmark id=24) ---
$bm = $d.Bookmarks("predictive-maintenance")
$rng = $bm.Range
$bm.Delete()
$d.Bookmarks.Add("automotive-safety", $rng)
$d.Paragraphs(8).Range.Text = "Automotive Safety"

# --- Body paragraph under "Automotive Safety" ---
$d.Paragraphs(9).Range.Text = @"
Newton’s laws are integral to automotive safety systems. Seat belts and airbags are designed based on the first law of motion, which states that an object in motion stays in motion unless acted upon by an external force. During a sudden stop or collision, seat belts restrain passengers to prevent harmful movement. Similarly, crumple zones in cars absorb and redirect forces according to Newton’s second law, reducing the impact on occupants. These applications have significantly reduced fatalities and injuries in vehicle accidents.
"@

# --- Heading3: "Fraud Detection" -> "Video Games" (bookmark id=25) ---
$bm = $d.Bookmarks("fraud-detection")
$rng = $bm.Range
$bm.Delete()
$d.Bookmarks.Add("video-games", $rng)
$d.Paragraphs(10).Range.Text = "Video Games"

# --- Body paragraph under "Video Games" ---
$d.Paragraphs(11).Range.Text = @"
Physics engines in video games use Newton’s laws to simulate realistic motion and interactions. The second law is used to calculate acceleration and deceleration of objects, while the third law ensures realistic collision effects. For example, in a racing game, the physics engine applies Newton’s laws to simulate how a car responds to braking, acceleration, and collisions. This creates a lifelike experience for players and enhances the overall gaming quality.
"@

# --- Heading3: "Text Classification" -> "Everyday Technology" (bookmark id=26) ---
$bm = $d.Bookmarks("text-classification")
$rng = $bm.Range
$bm.Delete()
$d.Bookmarks.Add("everyday-technology", $rng)
$d.Paragraphs(12).Range.Text = "Everyday Technology"

# --- Body paragraph under "Everyday Technology" ---
$d.Paragraphs(13).Range.Text = @"
Newton’s laws are applied in everyday technologies, such as washing machines and drying systems. During the spin cycle, clothes are subjected to motion, but water molecules remain at rest due to inertia (first law). The force exerted by the drum on the clothes (second law) removes water, and the interaction between the clothes and the drum follows the third law. These principles ensure efficient drying and prevent damage to fabrics.
"@

# --- Heading2: "Industry Case Study: Credit Card Fraud Detection" -> "Industry Case Study: SpaceX" (bookmark id=27) ---
$bm = $d.Bookmarks("X74c44ed3ad75ebfa9499fba50c218aa4884edf5")
$rng = $bm.Range
$bm.Delete()
$d.Bookmarks.Add("industry-case-study-spacex", $rng)
$d.Paragraphs(14).Range.Text = "Industry Case Study: SpaceX"

# --- Body paragraph under "Industry Case Study: SpaceX" ---
$d.Paragraphs(15).Range.Text = @"
SpaceX is a prime example of how Newton’s laws are applied in modern industry. The company uses Newton’s second law to calculate the thrust and acceleration of its rockets. For instance, the Falcon 9 rocket’s thrust-to-weight ratio is optimized based on the relationship F=ma (force equals mass times acceleration). Additionally, Newton’s third law is critical for understanding the reaction forces during propulsion, where expelled gases generate forward thrust. SpaceX’s reliance on these principles has enabled groundbreaking advancements in space exploration, including reusable rockets and precise orbital maneuvers.
"@

# --- Additional Resources hyperlinks: update link display text in place ---
# (scope Find to each list-item paragraph so only the hyperlink field run text
#  is replaced and the <w:hyperlink> field itself stays intact)
$p17 = $d.Paragraphs(17).Range
$p17.Find.Execute("Introduction to Support Vector Machines - IBM", $true, $false, $false, $false, $false, $true, 1, $false, "Newton’s Laws of Motion - Britannica", 2) | Out-Null

$p18 = $d.Paragraphs(18).Range
$p18.Find.Execute("Support Vector Machines Explained - Coursera", $true, $false, $false, $false, $false, $true, 1, $false, "Physics.org - Laws of Motion", 2) | Out-Null

$p19 = $d.Paragraphs(19).Range
$p19.Find.Execute("Recent Advances in SVM Research - Vilnius Tech", $true, $false, $false, $false, $false, $true, 1, $false, "Laws of Motion in Aerospace - Discover Engineering", 2) | Out-Null

$p20 = $d.Paragraphs(20).Range
$p20.Find.Execute("SVM for Anomaly Detection - GeeksforGeeks", $true, $false, $false, $false, $false, $true, 1, $false, "Robotics and Newton’s Laws - Ian McEachern", 2) | Out-Null

# --- Add a line break run right after the first three hyperlinks (4th stays as-is) ---
$d.Hyperlinks(1).Range.InsertAfter([char]11)
$d.Hyperlinks(2).Range.InsertAfter([char]11)
$d.Hyperlinks(3).Range.InsertAfter([char]11)
